$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.371.38'
$ws.Range('E2').Value = '  +1.85%  '
$ws.Range('D3').Value = '1.827.19'
$ws.Range('E3').Value = '  +2.90%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.49'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5353'
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4066'
$ws.Range('E8').Value = '  +9.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07597'
$ws.Range('E9').Value = '  +2.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.82'
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.104'
$ws.Range('E11').Value = '  +1.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.328'
$ws.Range('E12').Value = '  +4.46%  '
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.578'
$ws.Range('E14').Value = '  +5.51%  '
$ws.Range('E15').Value = '  +1.92%  '
$ws.Range('D16').Value = '1.825.86'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.30'
$ws.Range('E17').Value = '  +1.96%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001072'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06612'
$ws.Range('E19').Value = '  +3.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.61'
$ws.Range('E20').Value = '  +2.25%  '
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.080'
$ws.Range('E22').Value = '  +3.64%  '
$ws.Range('D23').Value = '28.395.29'
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.18'
$ws.Range('E24').Value = '  +1.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.183'
$ws.Range('E25').Value = '  +4.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.467'
$ws.Range('E26').Value = '  +8.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '157.77'
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.58'
$ws.Range('E28').Value = '  +2.36%  '
$ws.Range('D29').Value = '2.038.96'
$ws.Range('E29').Value = '  +3.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.63'
$ws.Range('E30').Value = '  +3.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.122'
$ws.Range('E31').Value = '  +2.06%  '
$ws.Range('E32').Value = '  +5.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.647'
$ws.Range('E33').Value = '  +3.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.640'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07242'
$ws.Range('E35').Value = '  +14.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2251'
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02335'
$ws.Range('E37').Value = '  +3.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.198'
$ws.Range('E38').Value = '  +5.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.809'
$ws.Range('E39').Value = '  +4.77%  '
$ws.Range('E40').Value = '  +2.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.28'
$ws.Range('E41').Value = '  +3.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.182'
$ws.Range('E42').Value = '  +1.28%  '
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.401'
$ws.Range('E44').Value = '  -1.99%  '
$ws.Range('E45').Value = '  +2.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.703'
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5846'
$ws.Range('E47').Value = '  +2.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.43'
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('E49').Value = '  +3.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.203'
$ws.Range('E50').Value = '  +2.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06887'
$ws.Range('E51').Value = '  +1.20%  '
